# Revert merge commit 237a1a
#
# Flip the "Featured" column (S) from "N" to "Y" for the four tool rows
# that were re-featured (Power Drill 3/8", 18 Gauge Brad Nailer,
# Sandblaster, Dwell Tachometer - rows 5-8), matching the Jig Saw row
# (row 9) whose "Featured" flag was already "Y".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S5").Value = "Y"
$ws.Range("S6").Value = "Y"
$ws.Range("S7").Value = "Y"
$ws.Range("S8").Value = "Y"

# Restore the author's on-screen selection/scroll position at save time.
[void]$ws.Range("S16").Select()
